# Update Sema3f-Nrp2 NATMI LR-pair values following Dr Hou advice
# (Ligand/Receptor-expressing cell counts changed 1 -> 3, with recalculated
# expression, specificity and edge-weight statistics for rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.44654466666667
$ws.Range("H2").Value = 58.339634
$ws.Range("I2").Value = 0.7934109702307454
$ws.Range("J2").Value = 0.7934109702307454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.52052333333333
$ws.Range("N2").Value = 103.56157
$ws.Range("O2").Value = 0.7684334662422598
$ws.Range("P2").Value = 0.7684334662422598
$ws.Range("Q2").Value = 671.3048989183756
$ws.Range("R2").Value = 6041.74409026538
$ws.Range("S2").Value = 0.6096835420090461
$ws.Range("T2").Value = 0.6096835420090461

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.44654466666667
$ws.Range("H3").Value = 58.339634
$ws.Range("I3").Value = 0.7934109702307454
$ws.Range("J3").Value = 0.7934109702307454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.347618333333333
$ws.Range("N3").Value = 16.042855
$ws.Range("O3").Value = 0.1190390091234806
$ws.Range("P3").Value = 0.1190390091234805
$ws.Range("Q3").Value = 103.9926987794522
$ws.Range("R3").Value = 935.93428901507
$ws.Range("S3").Value = 0.09444685572396726
$ws.Range("T3").Value = 0.09444685572396724

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.44654466666667
$ws.Range("H4").Value = 58.339634
$ws.Range("I4").Value = 0.7934109702307454
$ws.Range("J4").Value = 0.7934109702307454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.055101333333334
$ws.Range("N4").Value = 15.165304
$ws.Range("O4").Value = 0.1125275246342597
$ws.Range("P4").Value = 0.1125275246342597
$ws.Range("Q4").Value = 98.30425387319291
$ws.Range("R4").Value = 884.7382848587361
$ws.Range("S4").Value = 0.08928057249773209
$ws.Range("T4").Value = 0.08928057249773208

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.524415666666667
$ws.Range("H5").Value = 7.573247
$ws.Range("I5").Value = 0.1029951139231878
$ws.Range("J5").Value = 0.1029951139231878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.52052333333333
$ws.Range("N5").Value = 103.56157
$ws.Range("O5").Value = 0.7684334662422598
$ws.Range("P5").Value = 0.7684334662422598
$ws.Range("Q5").Value = 87.14414992419887
$ws.Range("R5").Value = 784.2973493177899
$ws.Range("S5").Value = 0.07914489239801165
$ws.Range("T5").Value = 0.07914489239801165

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.524415666666667
$ws.Range("H6").Value = 7.573247
$ws.Range("I6").Value = 0.1029951139231878
$ws.Range("J6").Value = 0.1029951139231878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.347618333333333
$ws.Range("N6").Value = 16.042855
$ws.Range("O6").Value = 0.1190390091234806
$ws.Range("P6").Value = 0.1190390091234805
$ws.Range("Q6").Value = 13.49961150002055
$ws.Range("R6").Value = 121.496503500185
$ws.Range("S6").Value = 0.01226043630597627
$ws.Range("T6").Value = 0.01226043630597627

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.524415666666667
$ws.Range("H7").Value = 7.573247
$ws.Range("I7").Value = 0.1029951139231878
$ws.Range("J7").Value = 0.1029951139231878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.055101333333334
$ws.Range("N7").Value = 15.165304
$ws.Range("O7").Value = 0.1125275246342597
$ws.Range("P7").Value = 0.1125275246342597
$ws.Range("Q7").Value = 12.76117700245422
$ws.Range("R7").Value = 114.850593022088
$ws.Range("S7").Value = 0.0115897852191999
$ws.Range("T7").Value = 0.0115897852191999

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.539092333333333
$ws.Range("H8").Value = 7.617277
$ws.Range("I8").Value = 0.1035939158460669
$ws.Range("J8").Value = 0.1035939158460669
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.52052333333333
$ws.Range("N8").Value = 103.56157
$ws.Range("O8").Value = 0.7684334662422598
$ws.Range("P8").Value = 0.7684334662422598
$ws.Range("Q8").Value = 87.6507961383211
$ws.Range("R8").Value = 788.8571652448899
$ws.Range("S8").Value = 0.07960503183520212
$ws.Range("T8").Value = 0.07960503183520212

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.539092333333333
$ws.Range("H9").Value = 7.617277
$ws.Range("I9").Value = 0.1035939158460669
$ws.Range("J9").Value = 0.1035939158460669
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.347618333333333
$ws.Range("N9").Value = 16.042855
$ws.Range("O9").Value = 0.1190390091234806
$ws.Range("P9").Value = 0.1190390091234805
$ws.Range("Q9").Value = 13.57809671175944
$ws.Range("R9").Value = 122.202870405835
$ws.Range("S9").Value = 0.01233171709353703
$ws.Range("T9").Value = 0.01233171709353703

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.539092333333333
$ws.Range("H10").Value = 7.617277
$ws.Range("I10").Value = 0.1035939158460669
$ws.Range("J10").Value = 0.1035939158460669
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.055101333333334
$ws.Range("N10").Value = 15.165304
$ws.Range("O10").Value = 0.1125275246342597
$ws.Range("P10").Value = 0.1125275246342597
$ws.Range("Q10").Value = 12.83536903968978
$ws.Range("R10").Value = 115.518321357208
$ws.Range("S10").Value = 0.01165716691732771
$ws.Range("T10").Value = 0.01165716691732771
